# Auto-generated Excel COM-interop script
# Applies updated market-price / profit values to the Chocobo Profits workbook
# as scraped by the scheduled runner (see commit message).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 49.66
$ws.Range("I15").Value = 49.66
$ws.Range("K15").Value = 148.98
$ws.Range("M15").Value = 20.02000000000001
$ws.Range("H88").Value = 5806.4165
$ws.Range("I88").Value = 4370
$ws.Range("J88").Value = 6524.625
$ws.Range("K88").Value = 4370
$ws.Range("L88").Value = 6524.625
$ws.Range("M88").Value = -3964
$ws.Range("N88").Value = -7336.625
$ws.Range("H91").Value = 5806.4165
$ws.Range("I91").Value = 4370
$ws.Range("J91").Value = 6524.625
$ws.Range("K91").Value = 4370
$ws.Range("L91").Value = 6524.625
$ws.Range("M91").Value = -2966
$ws.Range("N91").Value = -9332.625
$ws.Range("H137").Value = 672995.0600000001
$ws.Range("I137").Value = 2269124.2
$ws.Range("J137").Value = 2620.74
$ws.Range("K137").Value = 6807372.600000001
$ws.Range("L137").Value = 7862.219999999999
$ws.Range("M137").Value = -6804822.600000001
$ws.Range("N137").Value = -12962.22

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3032.56
$ws.Range("I32").Value = 2857.3635
$ws.Range("J32").Value = 3619.087
$ws.Range("K32").Value = 2857.3635
$ws.Range("L32").Value = 3619.087
$ws.Range("M32").Value = -2570.3635
$ws.Range("N32").Value = -4193.087
$ws.Range("H88").Value = 22224982
$ws.Range("I88").Value = 66666664
$ws.Range("J88").Value = 4140
$ws.Range("K88").Value = 66666664
$ws.Range("L88").Value = 4140
$ws.Range("M88").Value = -66666258
$ws.Range("N88").Value = -4952
$ws.Range("H91").Value = 22224982
$ws.Range("I91").Value = 66666664
$ws.Range("J91").Value = 4140
$ws.Range("K91").Value = 66666664
$ws.Range("L91").Value = 4140
$ws.Range("M91").Value = -66665260
$ws.Range("N91").Value = -6948
$ws.Range("H132").Value = 2471.12
$ws.Range("I132").Value = 1179.3846
$ws.Range("J132").Value = 3870.5
$ws.Range("K132").Value = 3538.1538
$ws.Range("L132").Value = 11611.5
$ws.Range("M132").Value = -1008.1538
$ws.Range("N132").Value = -16671.5
$ws.Range("H137").Value = 42158.332
$ws.Range("J137").Value = 42158.332
$ws.Range("L137").Value = 42158.332
$ws.Range("N137").Value = -52358.332

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 964.73334
$ws.Range("I64").Value = 961.6667
$ws.Range("J64").Value = 966.7778
$ws.Range("K64").Value = 961.6667
$ws.Range("L64").Value = 966.7778
$ws.Range("M64").Value = -736.6667
$ws.Range("N64").Value = -1416.7778
$ws.Range("H67").Value = 964.73334
$ws.Range("I67").Value = 961.6667
$ws.Range("J67").Value = 966.7778
$ws.Range("K67").Value = 961.6667
$ws.Range("L67").Value = 966.7778
$ws.Range("M67").Value = -181.6667
$ws.Range("N67").Value = -2526.7778
$ws.Range("H134").Value = 3414.353
$ws.Range("I134").Value = 1074.6428
$ws.Range("J134").Value = 6262.696
$ws.Range("K134").Value = 3223.9284
$ws.Range("L134").Value = 18788.088
$ws.Range("M134").Value = -688.9284000000002
$ws.Range("N134").Value = -23858.088
$ws.Range("H137").Value = 48750
$ws.Range("J137").Value = 48750
$ws.Range("L137").Value = 48750
$ws.Range("N137").Value = -58950

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4116782
$ws.Range("I16").Value = 7937898
$ws.Range("J16").Value = 1733.5385
$ws.Range("K16").Value = 7937898
$ws.Range("L16").Value = 1733.5385
$ws.Range("M16").Value = -7937611
$ws.Range("N16").Value = -2307.5385
$ws.Range("H31").Value = 194731.9
$ws.Range("I31").Value = 501408.44
$ws.Range("J31").Value = 2167.558
$ws.Range("K31").Value = 501408.44
$ws.Range("L31").Value = 2167.558
$ws.Range("M31").Value = -501113.44
$ws.Range("N31").Value = -2757.558
$ws.Range("H34").Value = 194731.9
$ws.Range("I34").Value = 501408.44
$ws.Range("J34").Value = 2167.558
$ws.Range("K34").Value = 501408.44
$ws.Range("L34").Value = 2167.558
$ws.Range("M34").Value = -501206.44
$ws.Range("N34").Value = -2571.558
$ws.Range("H58").Value = 2653.543
$ws.Range("I58").Value = 1498.5714
$ws.Range("J58").Value = 7273.4287
$ws.Range("K58").Value = 1498.5714
$ws.Range("L58").Value = 7273.4287
$ws.Range("M58").Value = -1295.5714
$ws.Range("N58").Value = -7679.4287
$ws.Range("H62").Value = 2925.625
$ws.Range("I62").Value = 3000.8333
$ws.Range("J62").Value = 2700
$ws.Range("K62").Value = 3000.8333
$ws.Range("L62").Value = 2700
$ws.Range("M62").Value = -2376.8333
$ws.Range("N62").Value = -3948
$ws.Range("H65").Value = 2925.625
$ws.Range("I65").Value = 3000.8333
$ws.Range("J65").Value = 2700
$ws.Range("K65").Value = 15004.1665
$ws.Range("L65").Value = 13500
$ws.Range("M65").Value = -11884.1665
$ws.Range("N65").Value = -19740
$ws.Range("H107").Value = 732.2632
$ws.Range("I107").Value = 410.53845
$ws.Range("J107").Value = 1429.3334
$ws.Range("K107").Value = 410.53845
$ws.Range("L107").Value = 1429.3334
$ws.Range("M107").Value = 1509.46155
$ws.Range("N107").Value = -5269.3334
$ws.Range("H113").Value = 4116782
$ws.Range("I113").Value = 7937898
$ws.Range("J113").Value = 1733.5385
$ws.Range("K113").Value = 7937898
$ws.Range("L113").Value = 1733.5385
$ws.Range("M113").Value = -7935728
$ws.Range("N113").Value = -6073.538500000001
$ws.Range("H132").Value = 2591.1538
$ws.Range("I132").Value = 2142.853
$ws.Range("J132").Value = 5639.6
$ws.Range("K132").Value = 6428.559
$ws.Range("L132").Value = 16918.8
$ws.Range("M132").Value = -3898.559
$ws.Range("N132").Value = -21978.8
$ws.Range("H134").Value = 1235.2727
$ws.Range("I134").Value = 814.5161000000001
$ws.Range("J134").Value = 7757
$ws.Range("K134").Value = 2443.5483
$ws.Range("L134").Value = 23271
$ws.Range("M134").Value = 91.45169999999962
$ws.Range("N134").Value = -28341
$ws.Range("H136").Value = 2653.543
$ws.Range("I136").Value = 1498.5714
$ws.Range("J136").Value = 7273.4287
$ws.Range("K136").Value = 4495.7142
$ws.Range("L136").Value = 21820.2861
$ws.Range("M136").Value = -1945.7142
$ws.Range("N136").Value = -26920.2861
$ws.Range("H139").Value = 33777.668
$ws.Range("J139").Value = 33777.668
$ws.Range("L139").Value = 33777.668
$ws.Range("N139").Value = -44057.668

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 478372.47
$ws.Range("I5").Value = 1252.1666
$ws.Range("K5").Value = 3756.4998
$ws.Range("M5").Value = -3644.4998
$ws.Range("H68").Value = 1594.6792
$ws.Range("I68").Value = 1349.5
$ws.Range("J68").Value = 1666.439
$ws.Range("K68").Value = 4048.5
$ws.Range("L68").Value = 4999.317
$ws.Range("M68").Value = -3237.5
$ws.Range("N68").Value = -6621.317
$ws.Range("H71").Value = 1594.6792
$ws.Range("I71").Value = 1349.5
$ws.Range("J71").Value = 1666.439
$ws.Range("K71").Value = 12145.5
$ws.Range("L71").Value = 14997.951
$ws.Range("M71").Value = -8089.5
$ws.Range("N71").Value = -23109.951
$ws.Range("H93").Value = 10240
$ws.Range("J93").Value = 10240
$ws.Range("L93").Value = 30720
$ws.Range("N93").Value = -34464
$ws.Range("H108").Value = 3668.4285
$ws.Range("I108").Value = 669.75
$ws.Range("K108").Value = 2009.25
$ws.Range("M108").Value = 870.75
$ws.Range("H126").Value = 4226
$ws.Range("I126").Value = 2710
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 8130
$ws.Range("L126").Value = 19500
$ws.Range("M126").Value = -3190
$ws.Range("N126").Value = -29380
$ws.Range("H135").Value = 478372.47
$ws.Range("I135").Value = 1252.1666
$ws.Range("K135").Value = 11269.4994
$ws.Range("M135").Value = -8734.499400000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13891877
$ws.Range("I80").Value = 22730380
$ws.Range("K80").Value = 22730380
$ws.Range("M80").Value = -22729382
$ws.Range("H83").Value = 13891877
$ws.Range("I83").Value = 22730380
$ws.Range("K83").Value = 113651900
$ws.Range("M83").Value = -113646908
$ws.Range("H126").Value = 3465.21
$ws.Range("I126").Value = 2809.3
$ws.Range("J126").Value = 4995.6665
$ws.Range("K126").Value = 8427.900000000001
$ws.Range("L126").Value = 14986.9995
$ws.Range("M126").Value = -5957.900000000001
$ws.Range("N126").Value = -19926.9995
$ws.Range("H132").Value = 3279.5715
$ws.Range("I132").Value = 1804.4
$ws.Range("K132").Value = 5413.200000000001
$ws.Range("M132").Value = -2883.200000000001
$ws.Range("H137").Value = 31870
$ws.Range("J137").Value = 43740
$ws.Range("L137").Value = 43740
$ws.Range("N137").Value = -53940

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1288.6774
$ws.Range("I82").Value = 802.76
$ws.Range("J82").Value = 3313.3333
$ws.Range("K82").Value = 802.76
$ws.Range("L82").Value = 3313.3333
$ws.Range("M82").Value = -441.76
$ws.Range("N82").Value = -4035.3333
$ws.Range("H85").Value = 1288.6774
$ws.Range("I85").Value = 802.76
$ws.Range("J85").Value = 3313.3333
$ws.Range("K85").Value = 802.76
$ws.Range("L85").Value = 3313.3333
$ws.Range("M85").Value = 445.24
$ws.Range("N85").Value = -5809.3333
$ws.Range("H132").Value = 3980.5588
$ws.Range("I132").Value = 3132.818
$ws.Range("J132").Value = 5534.75
$ws.Range("K132").Value = 9398.454000000002
$ws.Range("L132").Value = 16604.25
$ws.Range("M132").Value = -6868.454000000002
$ws.Range("N132").Value = -21664.25
$ws.Range("H136").Value = 2607.8276
$ws.Range("I136").Value = 1321.5428
$ws.Range("J136").Value = 4565.2173
$ws.Range("K136").Value = 3964.6284
$ws.Range("L136").Value = 13695.6519
$ws.Range("M136").Value = -1414.6284
$ws.Range("N136").Value = -18795.6519

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3821.2068
$ws.Range("I122").Value = 2369.2
$ws.Range("K122").Value = 7107.599999999999
$ws.Range("M122").Value = -4657.599999999999
$ws.Range("H132").Value = 11113544
$ws.Range("I132").Value = 1749.375
$ws.Range("J132").Value = 23812738
$ws.Range("K132").Value = 5248.125
$ws.Range("L132").Value = 71438214
$ws.Range("M132").Value = -2718.125
$ws.Range("N132").Value = -71443274
$ws.Range("H136").Value = 3536.258
$ws.Range("I136").Value = 1531.8823
$ws.Range("J136").Value = 5970.143
$ws.Range("K136").Value = 4595.6469
$ws.Range("L136").Value = 17910.429
$ws.Range("M136").Value = -2045.6469
$ws.Range("N136").Value = -23010.429
